# GitTreino.xlsx update — "Enviando a funcionalidade X"
# Adds new git-add / git-commit rows to the cheat-sheet table and
# tweaks the wording of the existing "git add" description.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 6: reword the "git add" explanation -----------------------------
$ws.Range("B6").Value = "Faz com que um arquivo seja enviado ao terminal e o autoriza a ser comitado e monitorado"

# --- Row 7: new entry, styled like row 6 (bold-ish header font, no fill) -
# Copy formatting only from A6 onto A7 so it keeps the same cell style
# (s="2") the rest of the "git add" family uses.
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B7").Value = "Adiciona todos os arquivos da pasta e subpasta ao terminal para ser comitado"
$ws.Range("A7").Value = "git add ."

# --- Row 8: new entry ------------------------------------------------------
$ws.Range("A8").Value = 'git commit (nome do arquivo) -m "Menssagem a ser salva"'
$ws.Range("B8").Value = "Commita um arquivo único e atrela a uma msg."

# --- Row 9: new entry ------------------------------------------------------
$ws.Range("A9").Value = 'git commit -a -m "Msg a registrar"'
$ws.Range("B9").Value = "Commita todos os arquivos de uma vez."

# --- Selection follows the new last row, like in the authored workbook ----
$ws.Range("A8").Select() | Out-Null
